# Generate Report for Handback
#
# The localization run for 8cace8e5-369e-4976-bf23-f08e7f682060.md has now
# come back "in sync with en-US" (handed back) instead of being stuck
# "Ready for handoff" / showing a stale-version error. Update the three
# report sheets (Overview, zh-cn, de-de) to reflect the new handback state.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
# Row 3 is the 8cace8e5-369e-4976-bf23-f08e7f682060.md file.
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 is the 8cace8e5-369e-4976-bf23-f08e7f682060.md file.
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-22 18:51:45"
$zhcn.Range("P3").Value = ""
$zhcn.Range("P1").EntireColumn.ColumnWidth = 12.8

# ---- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 3 is the 8cace8e5-369e-4976-bf23-f08e7f682060.md file.
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-22 18:51:53"
$dede.Range("P3").Value = ""
$dede.Range("P1").EntireColumn.ColumnWidth = 12.8
